# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values are updated for rows 2..74 (rows 35, 45, 63 stay 0 - unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 1;  3 = 1;  4 = 2;  5 = 1;  6 = 1;  7 = 2;  8 = 1;  9 = 2;  10 = 2;
    11 = 1; 12 = 0; 13 = 0; 14 = 3; 15 = 2; 16 = 3; 17 = 2; 18 = 2; 19 = 2; 20 = 2;
    21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 2; 26 = 1; 27 = 0; 28 = 3; 29 = 2; 30 = 2;
    31 = 2; 32 = 0; 33 = 2; 34 = 1; 36 = 0; 37 = 2; 38 = 2; 39 = 1; 40 = 1;
    41 = 3; 42 = 2; 43 = 1; 44 = 3; 46 = 2; 47 = 2; 48 = 1; 49 = 1; 50 = 3;
    51 = 2; 52 = 0; 53 = 1; 54 = 2; 55 = 0; 56 = 0; 57 = 0; 58 = 0; 59 = 2; 60 = 1;
    61 = 1; 62 = 2; 64 = 1; 65 = 2; 66 = 1; 67 = 2; 68 = 1; 69 = 0; 70 = 2;
    71 = 1; 72 = 1; 73 = 2; 74 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
